$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Origineel:" (Dutch, lang nl-NL).
# Its text becomes "Windows besturing:" and three brand-new paragraphs are
# inserted right after it: two new "Lijstalinea" bullet items documenting
# the extra PipeServer/PipeClient console steps, followed by a fresh
# (English, lang en-US) "Origineel:" paragraph that takes over the role the
# modified paragraph used to play.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd("`r`n`a") -eq "Origineel:") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Origineel:' paragraph"
}

# 1) Rewrite the existing run's text (keeps its original nl-NL formatting).
$target.Range.Text = "Windows besturing:"

# 2) Make room for the three new paragraphs right after it.
[void]$target.Range.InsertParagraphAfter()
[void]$target.Range.InsertParagraphAfter()
[void]$target.Range.InsertParagraphAfter()

$idx = $target.Index
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xmlService = "<w:p $ns><w:pPr><w:pStyle w:val=`"Lijstalinea`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"nl-NL`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"nl-NL`"/></w:rPr><w:t>Run Multitouch.Service.Console.exe</w:t></w:r></w:p>"
$xmlDriver  = "<w:p $ns><w:pPr><w:pStyle w:val=`"Lijstalinea`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:lang w:val=`"nl-NL`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"nl-NL`"/></w:rPr><w:t>Run Multitouch.Driver.Console.exe</w:t></w:r></w:p>"
$xmlOrig    = "<w:p $ns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Origineel:</w:t></w:r></w:p>"

# 3) Fill the three freshly-inserted empty paragraphs with their content.
[void]$d.Paragraphs($idx + 1).Range.InsertXML($xmlService)
[void]$d.Paragraphs($idx + 2).Range.InsertXML($xmlDriver)
[void]$d.Paragraphs($idx + 3).Range.InsertXML($xmlOrig)
